$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(
    -0.24739049268230673,
    -0.18271156586920068,
    -0.13110359287558282,
    -0.12247791454957735,
    -0.11890644397745476,
    -0.052743284399484835,
    -0.042124119963311113,
    -0.031968119192976374,
    -0.029674498310387598,
    -0.027427825424796382,
    -0.024412084622075447,
    -0.02088473140469338,
    -0.017179478783707403,
    -0.0090862807262475087,
    -0.0080560721897304788,
    -0.0060365806245790132,
    -0.0040034720231512111,
    -0.058005983811934669,
    -0.012091701022866719,
    -0.0080164041349135573,
    -0.0040056663984140428,
    -0.045719119303380573,
    -0.040504008245484613,
    -0.020099930939627697,
    -0.062205119436873346,
    -0.059664022587554655,
    -0.05711271102478932,
    -0.054817332261963614,
    -0.047626464170560823,
    0.01241955551196039,
    -0.029042953617045342,
    -0.01898179925493082
)

$bValues = @(
    0.24671754962511017,
    0.18082271212840872,
    0.13047791445741552,
    0.12190644392379468,
    0.11694840377778526,
    0.052124119833676374,
    0.041968119059516695,
    0.031674498244992133,
    0.029427825356340254,
    0.027412084544451432,
    0.024384731321569042,
    0.020679478697753773,
    0.017086280599897918,
    0.0090560721257340049,
    0.0080365805513356037,
    0.0060034719498696099,
    0.0039999999086690607,
    0.057863115314916058,
    0.012016404092952016,
    0.0080056663560679198,
    0.0039999999573279155,
    0.045504008183092637,
    0.040099930738168865,
    0.019999999795939694,
    0.062164022523104379,
    0.059612710958173665,
    0.05681733219033136,
    0.054626464047608181,
    0.047580443889033042,
    -0.012429897875388729,
    0.028981799102346884,
    0.018955321092457567
)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
